# Apply crypto price/volume updates to Sheet1 (columns D = Price, E = Volume(1h))
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.026.71"
$ws.Range("E2").Value = "  -0.58%  "

$ws.Range("D3").Value = "2.215.89"
$ws.Range("E3").Value = "  -1.48%  "

$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.67"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.96%  "

$ws.Range("E6").Value = "  +0.75%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.04"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.41%  "

$ws.Range("E8").Value = "  +0.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.607"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.56"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.41%  "

$ws.Range("E11").Value = "  +1.71%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.04"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.40%  "

$ws.Range("E13").Value = "  +0.67%  "

$ws.Range("D14").Value = "2.546.58"
$ws.Range("E14").Value = "  -1.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.28"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.838"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.09%  "

$ws.Range("D17").Value = "2.206.29"
$ws.Range("E17").Value = "  -1.49%  "

$ws.Range("D18").Value = "41.918.20"
$ws.Range("E18").Value = "  -0.62%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000108"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +7.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.19"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.72"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.51%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.78"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +17.40%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.59"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.53%  "

$ws.Range("E24").Value = "  -6.33%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.84"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.99%  "

$ws.Range("E26").Value = "  +0.09%  "

$ws.Range("E27").Value = "  +1.60%  "

$ws.Range("E28").Value = "  -1.31%  "

$ws.Range("E29").Value = "  -0.48%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.21"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.49"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.97%  "

$ws.Range("E32").Value = "  +9.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0795"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.34%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.96"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.78%  "

$ws.Range("E35").Value = "  -0.13%  "

$ws.Range("E36").Value = "  -10.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.27"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.09%  "

$ws.Range("E38").Value = "  -4.57%  "

$ws.Range("E39").Value = "  +0.56%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "65.52"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.79%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.12"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.00%  "

$ws.Range("E42").Value = "  -2.95%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.198"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.86%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.82"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.87%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "105.27"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.85%  "

$ws.Range("E46").Value = "  -1.67%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.42"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +5.05%  "

$ws.Range("E48").Value = "  -0.81%  "

$ws.Range("E49").Value = "  -0.71%  "

$ws.Range("E50").Value = "  -0.18%  "

$ws.Range("D51").Value = "2.423.77"
$ws.Range("E51").Value = "  -1.46%  "
